$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.768535614013672
$ws.Range("B1").Value = 2.810459613800049
$ws.Range("C1").Value = 2.630307197570801
$ws.Range("D1").Value = 3.472831964492798
$ws.Range("E1").Value = 5.136696338653564
